# Adding functionality to the Statistics Form
#
# Populates the "Transactions" sheet (headers: Id, Amount, Date, Type,
# Description already in row 1) with nine sample transaction rows, and
# refreshes the view/selection + forces a full recalc on load, matching
# what the Statistics Form now needs to read from this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("NMRUNEV7JO", "1233",   "08/04/2015", "RegularIncome",    "323231"),
    @("9IPFRR6AD1", "123",    "08/04/2015", "RegularIncome",    "123"),
    @("I43JUZDLPP", "2",      "08/04/2015", "RegularIncome",    "44"),
    @("S5Q1Z1H1Y1", "45.33",  "08/04/2015", "RegularIncome",    "asdlkjqaklqwjewe"),
    @("DZO3G7OHMZ", "44.32",  "08/04/2015", "IrregularExpense", "Food"),
    @("SVEHOIE5BM", "123.44", "08/04/2015", "RegularIncome",    "Salary"),
    @("Q9SPXOVFPD", "2193.3", "08/04/2015", "IrregularIncome",  "Jackpot"),
    @("ZWII13B4U3", "2213.3", "08/04/2015", "IrregularExpense", "Robbery"),
    @("CS521O7G9M", "233.3",  "08/04/2015", "RegularExpense",   "MobilePhone")
)

$cols = @("A", "B", "C", "D", "E")

$row = 2
foreach ($record in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $value = $record[$i]

        # Every value in this data set (ids, amounts, dates, types,
        # descriptions) is plain text. Amount/Date columns (and the
        # occasional numeric-looking Description) would otherwise be
        # auto-coerced into numbers/dates by Excel, so force a text
        # format on those cells before writing the value in.
        $numericLooking = $false
        if ($value -match '^[0-9]+(\.[0-9]+)?$') {
            $numericLooking = $true
        }
        if ($value -match '^\d{1,2}/\d{1,2}/\d{4}$') {
            $numericLooking = $true
        }
        if ($numericLooking) {
            $ws.Range($addr).NumberFormat = "@"
        }

        $ws.Range($addr).Value = $value
    }
    $row++
}

# Update the selection / view state to match the edited workbook.
$ws.DisplayRightToLeft = $false
$ws.Range("A3:XFD4").Select()

# Force a full recalculation on load.
$wb.ForceFullCalculation = $true
